# "Envio de mensajes Whatsapp"
# - Mark rows 26..65 (column T, the WhatsApp-sent flag) as "SI" (was "NO")
# - Append a new course row (row 70) for AG76-592 Fitopatologia General - 3259
# - Recompute the derived "items" summary counts (handled automatically by
#   the formulas already in place on the "items" sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Flip the "sent" flag (column T) from NO to SI for rows 26 through 65
# ---------------------------------------------------------------------
for ($r = 26; $r -le 65; $r++) {
    $ws.Range("T" + $r).Value = "SI"
}

# ---------------------------------------------------------------------
# 2) Append the new course row (row 70)
# ---------------------------------------------------------------------
$row = 70

$ws.Range("A$row").Value = 77
$ws.Range("B$row").Value = "AG"
$ws.Range("C$row").Value = "M"
$ws.Range("D$row").Value = "AG76 - 592 - FITOPATOLOGÍA GENERAL - 3259"
$ws.Range("E$row").Value = "VII"
$ws.Range("F$row").Value = "TARDE"
$ws.Range("G$row").Value = "B"
$ws.Range("H$row").Value = 20
$ws.Range("I$row").Value = "ORELLANA OZHO CARLOS MANUEL"
$ws.Range("S$row").Value = "https://chat.whatsapp.com/HCW8kplTG7DHTjzNrQJlKW"
$ws.Range("T$row").Value = "SI"
$ws.Range("U$row").Value = "https://aula.undc.edu.pe/course/view.php?id=375"

$ws.Range("V$row").Formula = "=MID(U$row,45,4)"
$ws.Range("W$row").Formula = "=MID(D$row,1,10)"
$ws.Range("X$row").Formula = "=TRIM(MID(D$row,14,222))"
$ws.Range("Y$row").Formula = "=TRIM(CONCATENATE(""AGRONOMIA "",E$row,""-"",F$row,""-"",G$row,"" "",LEFT(X$row,LEN(X$row)-7)))"
$ws.Range("Z$row").Formula = "=CONCATENATE(B$row,""_"",E$row,""-"",F$row,""-"",G$row,"" "",X$row)"
$ws.Range("AA$row").Formula = "=TRIM(MID(Z$row,1,25))"
$ws.Range("AB$row").Formula = "=CONCATENATE(""<p><a href='"",S$row,""' target='_blank'><img src='"",items!`$B`$1,""' alt='' width='291' height='42' role='presentation' class='img-responsive atto_image_button_text-bottom'></a><br></p>"")"

# ---------------------------------------------------------------------
# 3) Reflect the new selection left by the author while doing this edit
# ---------------------------------------------------------------------
$ws.Range("X30").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) Recalculate so every formula (including the items-sheet COUNTIFS /
#    SUM summary) carries a fresh cached value.
# ---------------------------------------------------------------------
$excel.CalculateFull() | Out-Null
